# Update "想去人数" (wish-to-attend count) values in the F column
# for rows 2, 3, 5, 6 on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 2344
    $ws.Range("F3").Value = 1832
    $ws.Range("F5").Value = 1128
    $ws.Range("F6").Value = 1062
}
